# Update Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 2-46, and
# append a new data row 47 (weekly refresh of the Oregano series).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = [DateTime]::FromOADate(44181)
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(2, 11).Value = 10000
$ws.Cells.Item(2, 12).Value = 12000
$ws.Cells.Item(2, 13).Value = 11000
$ws.Cells.Item(2, 16).Value = 3667

$ws.Cells.Item(3, 4).Value = [DateTime]::FromOADate(44377)
$ws.Cells.Item(3, 10).Value = 16
$ws.Cells.Item(3, 11).Value = 10000
$ws.Cells.Item(3, 12).Value = 10500
$ws.Cells.Item(3, 13).Value = 10250
$ws.Cells.Item(3, 16).Value = 3417

$ws.Cells.Item(4, 4).Value = [DateTime]::FromOADate(44566)
$ws.Cells.Item(4, 10).Value = 16
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11000
$ws.Cells.Item(4, 16).Value = 3667

$ws.Cells.Item(5, 4).Value = [DateTime]::FromOADate(44398)
$ws.Cells.Item(5, 10).Value = 16
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 10500
$ws.Cells.Item(5, 13).Value = 10250
$ws.Cells.Item(5, 16).Value = 3417

$ws.Cells.Item(6, 4).Value = [DateTime]::FromOADate(44314)
$ws.Cells.Item(6, 10).Value = 16
$ws.Cells.Item(6, 11).Value = 10000
$ws.Cells.Item(6, 12).Value = 10000
$ws.Cells.Item(6, 13).Value = 10000
$ws.Cells.Item(6, 16).Value = 3333

$ws.Cells.Item(7, 4).Value = [DateTime]::FromOADate(44447)
$ws.Cells.Item(7, 10).Value = 16
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 10500
$ws.Cells.Item(7, 13).Value = 10250
$ws.Cells.Item(7, 16).Value = 3417

$ws.Cells.Item(8, 4).Value = [DateTime]::FromOADate(44559)
$ws.Cells.Item(8, 10).Value = 7
$ws.Cells.Item(8, 11).Value = 10000
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = 11143
$ws.Cells.Item(8, 16).Value = 3714

$ws.Cells.Item(9, 4).Value = [DateTime]::FromOADate(44363)
$ws.Cells.Item(9, 10).Value = 16
$ws.Cells.Item(9, 11).Value = 10000
$ws.Cells.Item(9, 12).Value = 10000
$ws.Cells.Item(9, 13).Value = 10000
$ws.Cells.Item(9, 16).Value = 3333

$ws.Cells.Item(10, 4).Value = [DateTime]::FromOADate(44370)
$ws.Cells.Item(10, 10).Value = 16
$ws.Cells.Item(10, 11).Value = 10000
$ws.Cells.Item(10, 12).Value = 10500
$ws.Cells.Item(10, 13).Value = 10250
$ws.Cells.Item(10, 16).Value = 3417

$ws.Cells.Item(11, 4).Value = [DateTime]::FromOADate(44195)
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 10000
$ws.Cells.Item(11, 13).Value = 10000
$ws.Cells.Item(11, 16).Value = 3333

$ws.Cells.Item(12, 4).Value = [DateTime]::FromOADate(44468)
$ws.Cells.Item(12, 10).Value = 16
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 16).Value = 3500

$ws.Cells.Item(13, 4).Value = [DateTime]::FromOADate(44321)
$ws.Cells.Item(13, 10).Value = 25
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 10000
$ws.Cells.Item(13, 13).Value = 10000
$ws.Cells.Item(13, 16).Value = 3333

$ws.Cells.Item(14, 4).Value = [DateTime]::FromOADate(44342)
$ws.Cells.Item(14, 10).Value = 17
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = 10000
$ws.Cells.Item(14, 16).Value = 3333

$ws.Cells.Item(15, 4).Value = [DateTime]::FromOADate(44307)
$ws.Cells.Item(15, 10).Value = 160
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 16).Value = 3333

$ws.Cells.Item(16, 4).Value = [DateTime]::FromOADate(44435)
$ws.Cells.Item(16, 10).Value = 16
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 10500
$ws.Cells.Item(16, 13).Value = 10250
$ws.Cells.Item(16, 16).Value = 3417

$ws.Cells.Item(17, 4).Value = [DateTime]::FromOADate(44175)
$ws.Cells.Item(17, 10).Value = 70
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 12000
$ws.Cells.Item(17, 13).Value = 12000
$ws.Cells.Item(17, 16).Value = 4000

$ws.Cells.Item(18, 4).Value = [DateTime]::FromOADate(44419)
$ws.Cells.Item(18, 10).Value = 16
$ws.Cells.Item(18, 11).Value = 10000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = 10000
$ws.Cells.Item(18, 16).Value = 3333

$ws.Cells.Item(19, 4).Value = [DateTime]::FromOADate(44349)
$ws.Cells.Item(19, 10).Value = 12
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 10000
$ws.Cells.Item(19, 16).Value = 3333

$ws.Cells.Item(20, 4).Value = [DateTime]::FromOADate(44356)
$ws.Cells.Item(20, 10).Value = 16
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = 10000
$ws.Cells.Item(20, 16).Value = 3333

$ws.Cells.Item(21, 4).Value = [DateTime]::FromOADate(44540)
$ws.Cells.Item(21, 10).Value = 32
$ws.Cells.Item(21, 11).Value = 8500
$ws.Cells.Item(21, 12).Value = 9000
$ws.Cells.Item(21, 13).Value = 8719
$ws.Cells.Item(21, 16).Value = 2906

$ws.Cells.Item(22, 4).Value = [DateTime]::FromOADate(44391)
$ws.Cells.Item(22, 10).Value = 16
$ws.Cells.Item(22, 11).Value = 10000
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 10000
$ws.Cells.Item(22, 16).Value = 3333

$ws.Cells.Item(23, 4).Value = [DateTime]::FromOADate(44433)
$ws.Cells.Item(23, 10).Value = 16
$ws.Cells.Item(23, 11).Value = 10000
$ws.Cells.Item(23, 12).Value = 10500
$ws.Cells.Item(23, 13).Value = 10250
$ws.Cells.Item(23, 16).Value = 3417

$ws.Cells.Item(24, 4).Value = [DateTime]::FromOADate(44300)
$ws.Cells.Item(24, 10).Value = 16
$ws.Cells.Item(24, 11).Value = 10000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 10000
$ws.Cells.Item(24, 16).Value = 3333

$ws.Cells.Item(25, 4).Value = [DateTime]::FromOADate(44524)
$ws.Cells.Item(25, 10).Value = 16
$ws.Cells.Item(25, 11).Value = 9000
$ws.Cells.Item(25, 12).Value = 10000
$ws.Cells.Item(25, 13).Value = 9500
$ws.Cells.Item(25, 16).Value = 3167

$ws.Cells.Item(26, 4).Value = [DateTime]::FromOADate(44475)
$ws.Cells.Item(26, 10).Value = 16
$ws.Cells.Item(26, 11).Value = 9000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 9500
$ws.Cells.Item(26, 16).Value = 3167

$ws.Cells.Item(27, 4).Value = [DateTime]::FromOADate(44539)
$ws.Cells.Item(27, 10).Value = 16
$ws.Cells.Item(27, 11).Value = 9000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 9500
$ws.Cells.Item(27, 16).Value = 3167

$ws.Cells.Item(28, 4).Value = [DateTime]::FromOADate(44573)
$ws.Cells.Item(28, 10).Value = 16
$ws.Cells.Item(28, 11).Value = 10000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 11000
$ws.Cells.Item(28, 16).Value = 3667

$ws.Cells.Item(29, 4).Value = [DateTime]::FromOADate(44510)
$ws.Cells.Item(29, 10).Value = 16
$ws.Cells.Item(29, 11).Value = 9000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 13).Value = 9500
$ws.Cells.Item(29, 16).Value = 3167

$ws.Cells.Item(30, 4).Value = [DateTime]::FromOADate(44335)
$ws.Cells.Item(30, 10).Value = 16
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 16).Value = 3333

$ws.Cells.Item(31, 4).Value = [DateTime]::FromOADate(44272)
$ws.Cells.Item(31, 10).Value = 70
$ws.Cells.Item(31, 11).Value = 10000
$ws.Cells.Item(31, 12).Value = 10000
$ws.Cells.Item(31, 13).Value = 10000
$ws.Cells.Item(31, 16).Value = 3333

$ws.Cells.Item(32, 4).Value = [DateTime]::FromOADate(44489)
$ws.Cells.Item(32, 10).Value = 16
$ws.Cells.Item(32, 11).Value = 9000
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 9500
$ws.Cells.Item(32, 16).Value = 3167

$ws.Cells.Item(33, 4).Value = [DateTime]::FromOADate(44517)
$ws.Cells.Item(33, 10).Value = 16
$ws.Cells.Item(33, 11).Value = 9000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 9500
$ws.Cells.Item(33, 16).Value = 3167

$ws.Cells.Item(34, 4).Value = [DateTime]::FromOADate(44405)
$ws.Cells.Item(34, 10).Value = 16
$ws.Cells.Item(34, 11).Value = 10000
$ws.Cells.Item(34, 12).Value = 10500
$ws.Cells.Item(34, 13).Value = 10250
$ws.Cells.Item(34, 16).Value = 3417

$ws.Cells.Item(35, 4).Value = [DateTime]::FromOADate(44328)
$ws.Cells.Item(35, 10).Value = 16
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 12).Value = 10000
$ws.Cells.Item(35, 13).Value = 10000
$ws.Cells.Item(35, 16).Value = 3333

$ws.Cells.Item(36, 4).Value = [DateTime]::FromOADate(44552)
$ws.Cells.Item(36, 10).Value = 8
$ws.Cells.Item(36, 11).Value = 9000
$ws.Cells.Item(36, 12).Value = 10000
$ws.Cells.Item(36, 13).Value = 10000
$ws.Cells.Item(36, 16).Value = 3333

$ws.Cells.Item(37, 4).Value = [DateTime]::FromOADate(44426)
$ws.Cells.Item(37, 10).Value = 16
$ws.Cells.Item(37, 11).Value = 10000
$ws.Cells.Item(37, 12).Value = 10500
$ws.Cells.Item(37, 13).Value = 10250
$ws.Cells.Item(37, 16).Value = 3417

$ws.Cells.Item(38, 4).Value = [DateTime]::FromOADate(44482)
$ws.Cells.Item(38, 10).Value = 16
$ws.Cells.Item(38, 11).Value = 9000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 13).Value = 9500
$ws.Cells.Item(38, 16).Value = 3167

$ws.Cells.Item(39, 4).Value = [DateTime]::FromOADate(44461)
$ws.Cells.Item(39, 10).Value = 16
$ws.Cells.Item(39, 11).Value = 9500
$ws.Cells.Item(39, 12).Value = 10000
$ws.Cells.Item(39, 13).Value = 9750
$ws.Cells.Item(39, 16).Value = 3250

$ws.Cells.Item(40, 4).Value = [DateTime]::FromOADate(44545)
$ws.Cells.Item(40, 10).Value = 25
$ws.Cells.Item(40, 11).Value = 9000
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 13).Value = 9480
$ws.Cells.Item(40, 16).Value = 3160

$ws.Cells.Item(41, 4).Value = [DateTime]::FromOADate(44454)
$ws.Cells.Item(41, 10).Value = 16
$ws.Cells.Item(41, 11).Value = 9500
$ws.Cells.Item(41, 12).Value = 10000
$ws.Cells.Item(41, 13).Value = 9750
$ws.Cells.Item(41, 16).Value = 3250

$ws.Cells.Item(42, 4).Value = [DateTime]::FromOADate(44412)
$ws.Cells.Item(42, 10).Value = 25
$ws.Cells.Item(42, 11).Value = 10000
$ws.Cells.Item(42, 12).Value = 10500
$ws.Cells.Item(42, 13).Value = 10260
$ws.Cells.Item(42, 16).Value = 3420

$ws.Cells.Item(43, 4).Value = [DateTime]::FromOADate(44293)
$ws.Cells.Item(43, 10).Value = 16
$ws.Cells.Item(43, 11).Value = 10000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 13).Value = 10000
$ws.Cells.Item(43, 16).Value = 3333

$ws.Cells.Item(44, 4).Value = [DateTime]::FromOADate(44266)
$ws.Cells.Item(44, 10).Value = 160
$ws.Cells.Item(44, 11).Value = 10000
$ws.Cells.Item(44, 12).Value = 10000
$ws.Cells.Item(44, 13).Value = 10000
$ws.Cells.Item(44, 16).Value = 3333

$ws.Cells.Item(45, 4).Value = [DateTime]::FromOADate(44279)
$ws.Cells.Item(45, 10).Value = 16
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 13).Value = 10000
$ws.Cells.Item(45, 16).Value = 3333

$ws.Cells.Item(46, 4).Value = [DateTime]::FromOADate(44503)
$ws.Cells.Item(46, 10).Value = 16
$ws.Cells.Item(46, 11).Value = 8000
$ws.Cells.Item(46, 12).Value = 9000
$ws.Cells.Item(46, 13).Value = 8500
$ws.Cells.Item(46, 16).Value = 2833

# New row 47
$ws.Cells.Item(47, 1).Value = 9
$ws.Cells.Item(47, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(47, 3).Value = "Metropolitana"
$ws.Cells.Item(47, 4).Value = [DateTime]::FromOADate(44384)
$ws.Range("D47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 13
$ws.Cells.Item(47, 6).Value = 100112029
$ws.Cells.Item(47, 7).Value = "Orégano"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 25
$ws.Cells.Item(47, 11).Value = 10000
$ws.Cells.Item(47, 12).Value = 10500
$ws.Cells.Item(47, 13).Value = 10260
$ws.Cells.Item(47, 14).Value = "$/docena de atados"
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 3420
$ws.Cells.Item(47, 17).Value = 3
$ws.Cells.Item(47, 18).Value = "Hortaliza"
